# Update Dwayne Bravo's per-delivery stats: the two innings rows had their
# "runs"/"balls" figures transposed - row 2 should read 0/1 and row 3 should
# read 7/5. Source data is stored as text (t="str"), so write the new
# figures with a leading apostrophe to force text entry, then reset the
# cell style back to Normal so no stray number-format/quote-prefix styling
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'0"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "'1"
$ws.Range("D2").Style = "Normal"

$ws.Range("C3").Value = "'7"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "'5"
$ws.Range("D3").Style = "Normal"
